$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# E13: simple value correction (fuzzer bug fix)
$ws.Range("E13").Value = 332114255.60000002

# E15: simple value correction (fuzzer bug fix)
$ws.Range("E15").Value = 3000000

# E18: was a hard-coded value, restore it to a formula like its neighbors (D18, F18, G18)
$ws.Range("E18").Formula = "=SUM(E12:E17)"

# E21: was a hard-coded value, restore it to a formula like its neighbors (D21, F21, G21)
$ws.Range("E21").Formula = "=SUM(E18:E20)"

# E26: value correction
$ws.Range("E26").Value = 1018613404

# Recalculate so dependent formula cells (E23, E25, E28, F28, C29, etc.) update
$excel.CalculateFullRebuild()
